# Rewrite sheet8 ("tbl8") into the expanded "Cumulative number strategies" table
# with per-threshold (0.1/0.2/0.5/0.6/0.7/0.8/1) coefficient + p-value columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(8)

# Header row (row 1)
$ws.Range("A1").Value = "Cumulative number strategies"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "0.1_Coefficient (95% CI)"
$ws.Range("D1").Value = "0.1_p-value"
$ws.Range("E1").Value = "0.2_Coefficient (95% CI)"
$ws.Range("F1").Value = "0.2_p-value"
$ws.Range("G1").Value = "0.5_Coefficient (95% CI)"
$ws.Range("H1").Value = "0.5_p-value"
$ws.Range("I1").Value = "0.6_Coefficient (95% CI)"
$ws.Range("J1").Value = "0.6_p-value"
$ws.Range("K1").Value = "0.7_Coefficient (95% CI)"
$ws.Range("L1").Value = "0.7_p-value"
$ws.Range("M1").Value = "0.8_Coefficient (95% CI)"
$ws.Range("N1").Value = "0.8_p-value"
$ws.Range("O1").Value = "1_Coefficient (95% CI)"
$ws.Range("P1").Value = "1_p-value"

# Data rows
$ws.Range("A2").Value = "Intercept"
$ws.Range("B2").Value = "(Intercept)"
$ws.Range("C2").Value = "1.02 (0.31, 1.74)"
$ws.Range("D2").Value = 0.05
$ws.Range("E2").Value = "0.97 (0.25, 1.73)"
$ws.Range("F2").Value = 0.06
$ws.Range("G2").Value = "0.98 (0.21, 1.65)"
$ws.Range("H2").Value = 0.06
$ws.Range("I2").Value = "0.94 (0.18, 1.62)"
$ws.Range("J2").Value = 0.06
$ws.Range("K2").Value = "0.92 (0.21, 1.66)"
$ws.Range("L2").Value = 0.07000000000000001
$ws.Range("M2").Value = "0.94 (0.14, 1.76)"
$ws.Range("N2").Value = 0.07000000000000001
$ws.Range("O2").Value = "0.96 (0.13, 1.73)"
$ws.Range("P2").Value = 0.07000000000000001
$ws.Range("A3").Value = "1"
$ws.Range("B3").Value = "ss1"
$ws.Range("C3").Value = "0.15 (-0.53, 0.76)"
$ws.Range("D3").Value = 0.6
$ws.Range("E3").Value = "1.06 (0.34, 1.77)"
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = "1.28 (0.43, 2.21)"
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = "1.22 (0.42, 2.05)"
$ws.Range("J3").Value = 0.01
$ws.Range("K3").Value = "1.23 (0.37, 2.22)"
$ws.Range("L3").Value = 0.01
$ws.Range("M3").Value = "1.25 (0.38, 2.21)"
$ws.Range("N3").Value = 0.01
$ws.Range("O3").Value = "1.16 (0.2, 2.21)"
$ws.Range("P3").Value = 0.03
$ws.Range("A4").Value = "2"
$ws.Range("B4").Value = "ss2"
$ws.Range("C4").Value = "-0.23 (-0.77, 0.31)"
$ws.Range("D4").Value = 0.39
$ws.Range("E4").Value = "-0.2 (-0.85, 0.37)"
$ws.Range("F4").Value = 0.52
$ws.Range("G4").Value = "-0.04 (-0.72, 0.62)"
$ws.Range("H4").Value = 0.92
$ws.Range("I4").Value = "-0.07 (-0.83, 0.78)"
$ws.Range("J4").Value = 0.86
$ws.Range("K4").Value = "-0.26 (-1.12, 0.57)"
$ws.Range("L4").Value = 0.53
$ws.Range("M4").Value = "0.01 (-0.8, 0.81)"
$ws.Range("N4").Value = 0.98
$ws.Range("O4").Value = "0.13 (-0.63, 0.98)"
$ws.Range("P4").Value = 0.75
$ws.Range("A5").Value = "3"
$ws.Range("B5").Value = "ss3"
$ws.Range("C5").Value = "-0.53 (-1.05, -0.02)"
$ws.Range("D5").Value = 0.04
$ws.Range("E5").Value = "-0.12 (-0.68, 0.39)"
$ws.Range("F5").Value = 0.68
$ws.Range("G5").Value = "-0.22 (-0.78, 0.37)"
$ws.Range("H5").Value = 0.45
$ws.Range("I5").Value = "0.12 (-0.5, 0.73)"
$ws.Range("J5").Value = 0.71
$ws.Range("K5").Value = "0.01 (-0.75, 0.68)"
$ws.Range("L5").Value = 0.99
$ws.Range("M5").Value = "-0.26 (-1.02, 0.5)"
$ws.Range("N5").Value = 0.49
$ws.Range("O5").Value = "-0.24 (-1.11, 0.52)"
$ws.Range("P5").Value = 0.53
$ws.Range("A6").Value = "Percent two or more races"
$ws.Range("B6").Value = "percenttwoormoreraces"
$ws.Range("C6").Value = "0.24 (0.03, 0.46)"
$ws.Range("D6").Value = 0.03
$ws.Range("E6").Value = "0.24 (0, 0.46)"
$ws.Range("F6").Value = 0.03
$ws.Range("G6").Value = "0.25 (0.06, 0.48)"
$ws.Range("H6").Value = 0.02
$ws.Range("I6").Value = "0.25 (0.06, 0.45)"
$ws.Range("J6").Value = 0.03
$ws.Range("K6").Value = "0.25 (0.06, 0.46)"
$ws.Range("L6").Value = 0.03
$ws.Range("M6").Value = "0.26 (0.03, 0.49)"
$ws.Range("N6").Value = 0.02
$ws.Range("O6").Value = "0.26 (0.05, 0.47)"
$ws.Range("P6").Value = 0.02
$ws.Range("A7").Value = "Percent Asian"
$ws.Range("B7").Value = "percentasian"
$ws.Range("C7").Value = "0.01 (-0.2, 0.22)"
$ws.Range("D7").Value = 0.9399999999999999
$ws.Range("E7").Value = "0.04 (-0.16, 0.27)"
$ws.Range("F7").Value = 0.68
$ws.Range("G7").Value = "0.04 (-0.19, 0.23)"
$ws.Range("H7").Value = 0.7
$ws.Range("I7").Value = "0.03 (-0.18, 0.23)"
$ws.Range("J7").Value = 0.76
$ws.Range("K7").Value = "0.05 (-0.17, 0.25)"
$ws.Range("L7").Value = 0.65
$ws.Range("M7").Value = "0.05 (-0.14, 0.26)"
$ws.Range("N7").Value = 0.63
$ws.Range("O7").Value = "0.05 (-0.14, 0.24)"
$ws.Range("P7").Value = 0.63
$ws.Range("A8").Value = "Percent White"
$ws.Range("B8").Value = "percentwhite"
$ws.Range("C8").Value = "-0.19 (-0.57, 0.15)"
$ws.Range("D8").Value = 0.29
$ws.Range("E8").Value = "-0.16 (-0.51, 0.23)"
$ws.Range("F8").Value = 0.39
$ws.Range("G8").Value = "-0.16 (-0.57, 0.22)"
$ws.Range("H8").Value = 0.38
$ws.Range("I8").Value = "-0.19 (-0.56, 0.19)"
$ws.Range("J8").Value = 0.3
$ws.Range("K8").Value = "-0.15 (-0.49, 0.22)"
$ws.Range("L8").Value = 0.43
$ws.Range("M8").Value = "-0.16 (-0.55, 0.19)"
$ws.Range("N8").Value = 0.37
$ws.Range("O8").Value = "-0.17 (-0.51, 0.19)"
$ws.Range("P8").Value = 0.35
$ws.Range("A9").Value = "Percent free and reduced lunch"
$ws.Range("B9").Value = "percentfreereducedlunch"
$ws.Range("C9").Value = "-0.28 (-0.57, 0.06)"
$ws.Range("D9").Value = 0.08
$ws.Range("E9").Value = "-0.22 (-0.54, 0.08)"
$ws.Range("F9").Value = 0.17
$ws.Range("G9").Value = "-0.23 (-0.55, 0.09)"
$ws.Range("H9").Value = 0.14
$ws.Range("I9").Value = "-0.26 (-0.57, 0.06)"
$ws.Range("J9").Value = 0.1
$ws.Range("K9").Value = "-0.24 (-0.52, 0.1)"
$ws.Range("L9").Value = 0.14
$ws.Range("M9").Value = "-0.23 (-0.55, 0.06)"
$ws.Range("N9").Value = 0.14
$ws.Range("O9").Value = "-0.24 (-0.54, 0.07)"
$ws.Range("P9").Value = 0.14
$ws.Range("A10").Value = "SVI Overall Rank"
$ws.Range("B10").Value = "rplthemes"
$ws.Range("C10").Value = "0.03 (-0.23, 0.25)"
$ws.Range("D10").Value = 0.82
$ws.Range("E10").Value = "0 (-0.23, 0.25)"
$ws.Range("F10").Value = 0.99
$ws.Range("G10").Value = "0.02 (-0.21, 0.26)"
$ws.Range("H10").Value = 0.9
$ws.Range("I10").Value = "0.02 (-0.25, 0.27)"
$ws.Range("J10").Value = 0.84
$ws.Range("K10").Value = "0.04 (-0.21, 0.3)"
$ws.Range("L10").Value = 0.76
$ws.Range("M10").Value = "0.05 (-0.18, 0.28)"
$ws.Range("N10").Value = 0.66
$ws.Range("O10").Value = "0.02 (-0.19, 0.25)"
$ws.Range("P10").Value = 0.85
$ws.Range("A11").Value = "Percent Black or African American"
$ws.Range("B11").Value = "percentblackorafricanamerican"
$ws.Range("C11").Value = "-0.07 (-0.34, 0.17)"
$ws.Range("D11").Value = 0.5600000000000001
$ws.Range("E11").Value = "-0.03 (-0.29, 0.21)"
$ws.Range("F11").Value = 0.79
$ws.Range("G11").Value = "-0.05 (-0.32, 0.22)"
$ws.Range("H11").Value = 0.6899999999999999
$ws.Range("I11").Value = "-0.04 (-0.27, 0.22)"
$ws.Range("J11").Value = 0.74
$ws.Range("K11").Value = "-0.03 (-0.28, 0.24)"
$ws.Range("L11").Value = 0.85
$ws.Range("M11").Value = "-0.07 (-0.31, 0.16)"
$ws.Range("N11").Value = 0.58
$ws.Range("O11").Value = "-0.06 (-0.32, 0.22)"
$ws.Range("P11").Value = 0.66
$ws.Range("A12").Value = "High school"
$ws.Range("B12").Value = "schoollevelHS"
$ws.Range("C12").Value = "0.28 (-0.18, 0.75)"
$ws.Range("D12").Value = 0.24
$ws.Range("E12").Value = "0.23 (-0.25, 0.69)"
$ws.Range("F12").Value = 0.33
$ws.Range("G12").Value = "0.29 (-0.22, 0.79)"
$ws.Range("H12").Value = 0.23
$ws.Range("I12").Value = "0.29 (-0.17, 0.78)"
$ws.Range("J12").Value = 0.24
$ws.Range("K12").Value = "0.35 (-0.13, 0.79)"
$ws.Range("L12").Value = 0.15
$ws.Range("M12").Value = "0.39 (-0.14, 0.88)"
$ws.Range("N12").Value = 0.11
$ws.Range("O12").Value = "0.33 (-0.14, 0.83)"
$ws.Range("P12").Value = 0.19
$ws.Range("A13").Value = "Middle school"
$ws.Range("B13").Value = "schoollevelMS"
$ws.Range("C13").Value = "0.24 (-0.24, 0.65)"
$ws.Range("D13").Value = 0.31
$ws.Range("E13").Value = "0.19 (-0.23, 0.61)"
$ws.Range("F13").Value = 0.41
$ws.Range("G13").Value = "0.27 (-0.16, 0.72)"
$ws.Range("H13").Value = 0.25
$ws.Range("I13").Value = "0.25 (-0.22, 0.71)"
$ws.Range("J13").Value = 0.3
$ws.Range("K13").Value = "0.31 (-0.12, 0.79)"
$ws.Range("L13").Value = 0.19
$ws.Range("M13").Value = "0.28 (-0.18, 0.73)"
$ws.Range("N13").Value = 0.23
$ws.Range("O13").Value = "0.25 (-0.23, 0.74)"
$ws.Range("P13").Value = 0.3
$ws.Range("A14").Value = "Change in county COVID-19 case rate"
$ws.Range("B14").Value = "cntycaseschange"
$ws.Range("C14").Value = "-0.09 (-0.47, 0.22)"
$ws.Range("D14").Value = 0.61
$ws.Range("E14").Value = "-0.13 (-0.49, 0.2)"
$ws.Range("F14").Value = 0.43
$ws.Range("G14").Value = "-0.12 (-0.44, 0.19)"
$ws.Range("H14").Value = 0.49
$ws.Range("I14").Value = "-0.08 (-0.41, 0.23)"
$ws.Range("J14").Value = 0.63
$ws.Range("K14").Value = "-0.08 (-0.42, 0.25)"
$ws.Range("L14").Value = 0.62
$ws.Range("M14").Value = "-0.06 (-0.37, 0.25)"
$ws.Range("N14").Value = 0.7
$ws.Range("O14").Value = "-0.07 (-0.4, 0.29)"
$ws.Range("P14").Value = 0.67
$ws.Range("A15").Value = "4"
$ws.Range("B15").Value = "ss4"
$ws.Range("E15").Value = "-0.24 (-0.85, 0.3)"
$ws.Range("F15").Value = 0.41
$ws.Range("G15").Value = "-0.28 (-0.86, 0.3)"
$ws.Range("H15").Value = 0.37
$ws.Range("I15").Value = "-0.25 (-0.8, 0.34)"
$ws.Range("J15").Value = 0.41
$ws.Range("K15").Value = "0.1 (-0.51, 0.81)"
$ws.Range("L15").Value = 0.76
$ws.Range("M15").Value = "0.46 (-0.38, 1.17)"
$ws.Range("N15").Value = 0.24
$ws.Range("O15").Value = "0.46 (-0.29, 1.21)"
$ws.Range("P15").Value = 0.24
$ws.Range("A16").Value = "5"
$ws.Range("B16").Value = "ss5"
$ws.Range("E16").Value = "-0.84 (-1.63, -0.15)"
$ws.Range("F16").Value = 0.03
$ws.Range("G16").Value = "-0.1 (-0.78, 0.56)"
$ws.Range("H16").Value = 0.77
$ws.Range("I16").Value = "-0.19 (-0.85, 0.5)"
$ws.Range("J16").Value = 0.6
$ws.Range("K16").Value = "-0.26 (-0.94, 0.4)"
$ws.Range("L16").Value = 0.44
$ws.Range("M16").Value = "-0.26 (-0.89, 0.36)"
$ws.Range("N16").Value = 0.38
$ws.Range("O16").Value = "-0.19 (-0.83, 0.41)"
$ws.Range("P16").Value = 0.55
$ws.Range("A17").Value = "6"
$ws.Range("B17").Value = "ss6"
$ws.Range("G17").Value = "-0.93 (-1.89, -0.12)"
$ws.Range("H17").Value = 0.03
$ws.Range("I17").Value = "-0.4 (-1.13, 0.22)"
$ws.Range("J17").Value = 0.26
$ws.Range("K17").Value = "-0.38 (-1.08, 0.27)"
$ws.Range("L17").Value = 0.27
$ws.Range("M17").Value = "-0.67 (-1.4, 0.02)"
$ws.Range("N17").Value = 0.06
$ws.Range("O17").Value = "-0.72 (-1.45, 0.04)"
$ws.Range("P17").Value = 0.05
$ws.Range("A18").Value = "7"
$ws.Range("B18").Value = "ss7"
$ws.Range("I18").Value = "-0.76 (-1.78, 0.29)"
$ws.Range("J18").Value = 0.15
$ws.Range("K18").Value = "-0.21 (-1.2, 0.69)"
$ws.Range("L18").Value = 0.64
$ws.Range("M18").Value = "0.23 (-0.5, 0.93)"
$ws.Range("N18").Value = 0.57
$ws.Range("O18").Value = "-0.12 (-0.86, 0.69)"
$ws.Range("P18").Value = 0.75
$ws.Range("A19").Value = "8"
$ws.Range("B19").Value = "ss8"
$ws.Range("K19").Value = "-0.78 (-1.73, 0.21)"
$ws.Range("L19").Value = 0.14
$ws.Range("M19").Value = "-0.82 (-1.78, 0.11)"
$ws.Range("N19").Value = 0.09
$ws.Range("O19").Value = "-0.41 (-1.25, 0.49)"
$ws.Range("P19").Value = 0.35
$ws.Range("A20").Value = "9"
$ws.Range("B20").Value = "ss9"
$ws.Range("M20").Value = "-0.4 (-2.34, 1.49)"
$ws.Range("N20").Value = 0.68
$ws.Range("O20").Value = "-0.46 (-2.44, 1.35)"
$ws.Range("P20").Value = 0.63
